# Update the "Fitness" log values (column C) on the active sheet.
# Rows 2-60  : Fitness 7534 / 7345 / 7312 / 7310 -> 7295
# Rows 61-93 : Fitness 7310 / 7295             -> 7293
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C60").Value = 7295
$ws.Range("C61:C93").Value = 7293
